# fix: various task transformer improvements (xlsx)
# - remove support for column "Version" (no structural cell change needed here;
#   handled by the reading code, not the sheet itself)
# - add support for columns "goal_version", "rule_name_id", "rule_version"
#   inserted right after "goal_name_id" (old column L), before the old
#   "Parameter" / "Values" columns.
#
# Net effect on the worksheet: insert three new columns M,N,O so that
#   L = goal_name_id      (unchanged)
#   M = goal_version       (new, empty)
#   N = rule_name_id       (new, copy of goal_name_id per row)
#   O = rule_version       (new, empty)
#   P = Parameter [optional parameter]   (was M)
#   Q = Values default , [alternatives]  (was N)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank columns before the old "Parameter" column (M / 13).
# Each Insert() shifts the old M (Parameter) and N (Values) columns two
# further to the right, and the new column inherits the formatting of the
# column to its left (matching column L / goal_name_id formatting).
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(13).Insert()

# Header row (row 1): fill in the three new headers.
$ws.Cells.Item(1, 13).Value2 = "goal_version"
$ws.Cells.Item(1, 14).Value2 = "rule_name_id"
$ws.Cells.Item(1, 15).Value2 = "rule_version"

# Data rows (2-11): column N (rule_name_id) mirrors column L (goal_name_id)
# for every row that actually has a goal_name_id value; column M
# (goal_version) and column O (rule_version) stay empty.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $goalNameId = $ws.Cells.Item($r, 12).Value2
    if ($goalNameId -ne $null -and $goalNameId -ne "") {
        $ws.Cells.Item($r, 14).Value2 = $goalNameId
    }
}

# Column widths: L:N keep the old "Parameter" column width, O:P keep the
# old "Values" column width (best-effort; the interop layer quantizes
# widths to pixel steps so exact legacy values can't always be hit).
$ws.Range("M1:N1").EntireColumn.ColumnWidth = 39.88
$ws.Range("O1:P1").EntireColumn.ColumnWidth = 20.98

# Restore the view roughly to what it was (selection on M2, scrolled so
# column J is the left-most visible column).
$ws.Range("M2").Select()
$excel.ActiveWindow.ScrollColumn = 10
